$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column K (최종점수 / final score) for rows 2-7
$ws.Range("K2").Value = 56.6
$ws.Range("K3").Value = 56.6
$ws.Range("K4").Value = 52.4
$ws.Range("K5").Value = 47.4
$ws.Range("K6").Value = 45.2
$ws.Range("K7").Value = 44.6

# Update column N (MACRO_SCORE) for rows 2-7
$ws.Range("N2:N7").Value = 54.02451352198364
